$wb = $excel.ActiveWorkbook

# --- Sheet "Diario-Restante": task "Configurar la carga de modelos desde GitHub" (row 10)
#     and "Scrum Team" row (row 11) get effort numbers filled in ---
$wsRestante = $wb.Worksheets.Item("Diario-Restante")
$wsRestante.Range("F10").Value = 1
$wsRestante.Range("C11").Value = "Scrum Team"
$wsRestante.Range("D11").Value = 3
$wsRestante.Range("E11").Value = 2

# --- Sheet "Diario-Realizado": same task row (row 10) gets its effort filled in ---
$wsRealizado = $wb.Worksheets.Item("Diario-Realizado")
$wsRealizado.Range("C10").Value = "Scrum Team"
$wsRealizado.Range("D10").Value = 3
$wsRealizado.Range("E10").Value = 2
$wsRealizado.Range("F10").Value = 1

$wb.Application.CalculateFullRebuild()
